$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WorstCell Threshold")
$ws.Activate()

# Copy formatting from row 8 into row 9 first, to inherit borders/fonts/fills/row height
$ws.Range("A8:E8").Copy()
$ws.Range("A9:E9").PasteSpecial(-4122) # xlPasteFormats
$ws.Application.CutCopyMode = $false

# B9 inherits the A-column look (same font/fill/border as A9) but centered
$ws.Range("A8").Copy()
$ws.Range("B9").PasteSpecial(-4122) # xlPasteFormats
$ws.Application.CutCopyMode = $false
$ws.Range("B9").HorizontalAlignment = -4108 # xlCenter

$ws.Range("A9").Value = "5G Total Traffic Volume"
$ws.Range("B9").Value = "NR_50823a"
$ws.Range("C9").Value = "[GB]"
$ws.Range("D9").Value = "<"
$ws.Range("E9").Value = 100

$ws.Rows.Item(9).RowHeight = 15.75

$ws.Range("E10").Select() | Out-Null
